$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.016.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.095.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.00'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.68%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.093.32'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.37'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.613.85'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.989.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.02'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.096.97'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.21'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.686'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.40'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.10'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0941'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '47.53'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.71%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.942'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.313'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.02'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.08'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.799.22'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.76'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.53'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.86%  '

# Row 43/44 swap: dogwifhat <-> Cosmos (with updated values)
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.71"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.66%  "

# Row 46/47 swap: Bittensor <-> VeChain (with updated values)
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0345"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "368.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.09%  "
